# Add links to PPT for Custom Core
#
# Splits three text runs so that a sub-string of each becomes its own run
# with a hyperlink (a:hlinkClick) applied via ActionSettings, matching the
# target OOXML diff:
#   - Slide 10: "You can download the source files here. "
#       -> "...source files " + "here" (hyperlink) + ". "
#   - Slide 12: "Here is the source code for the modified "
#       -> "Here" (hyperlink) + " is the source code for the modified "
#   - Slide 3:  "Enter the Verilog code to add two 8-bit numbers"
#       -> "Enter the " + "Verilog code" (hyperlink) + " to add two 8-bit numbers"

$p = $ppt.ActivePresentation

function Add-RunHyperlink($Slide, $ShapeIndex, $ParagraphIndex, $SubStart, $SubLength, $Url) {
    $shape = $Slide.Shapes.Item($ShapeIndex)
    $textRange = $shape.TextFrame.TextRange
    $paragraph = $textRange.Paragraphs($ParagraphIndex)

    $absoluteStart = $paragraph.Start + $SubStart - 1
    $linkRange = $textRange.Characters($absoluteStart, $SubLength)
    $linkRange.ActionSettings.Item(1).Hyperlink.Address = $Url
}

$sourceUrl = "https://github.com/Future-Electronics-Design-Center/Creative-Eval-Board"

# Slide 10 ("You can download the source files here. ") -> link "here"
$slide10 = $p.Slides.Item(10)
Add-RunHyperlink $slide10 2 7 35 4 $sourceUrl

# Slide 12 ("Here is the source code for the modified main.c") -> link "Here"
$slide12 = $p.Slides.Item(12)
Add-RunHyperlink $slide12 2 5 1 4 $sourceUrl

# Slide 3 ("Enter the Verilog code to add two 8-bit numbers") -> link "Verilog code"
$slide3 = $p.Slides.Item(3)
Add-RunHyperlink $slide3 2 1 11 12 $sourceUrl

Write-Output "Hyperlinks added on slides 3, 10, 12"
